$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "62.494.83"
$ws.Range("E2").Value = "  -6.26%  "

# Row 3
Set-TextValue "D3" "3.097.40"
$ws.Range("E3").Value = "  -6.52%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
Set-TextValue "D5" "556.48"
$ws.Range("E5").Value = "  -5.96%  "

# Row 6
Set-TextValue "D6" "160.47"
$ws.Range("E6").Value = "  -11.59%  "

# Row 7
$ws.Range("E7").Value = "  +0.12%  "

# Row 8
Set-TextValue "D8" "0.578"
$ws.Range("E8").Value = "  -9.39%  "

# Row 9
Set-TextValue "D9" "3.091.79"
$ws.Range("E9").Value = "  -6.58%  "

# Row 10
Set-TextValue "D10" "6.70"
$ws.Range("E10").Value = "  -2.52%  "

# Row 11
Set-TextValue "D11" "0.114"
$ws.Range("E11").Value = "  -9.51%  "

# Row 12
Set-TextValue "D12" "0.375"
$ws.Range("E12").Value = "  -6.99%  "

# Row 13
Set-TextValue "D13" "3.639.67"
$ws.Range("E13").Value = "  -6.41%  "

# Row 14
$ws.Range("E14").Value = "  -2.08%  "

# Row 15
Set-TextValue "D15" "62.595.93"
$ws.Range("E15").Value = "  -6.08%  "

# Row 16
Set-TextValue "D16" "24.42"
$ws.Range("E16").Value = "  -8.71%  "

# Row 17
Set-TextValue "D17" "3.093.83"
$ws.Range("E17").Value = "  -6.19%  "

# Row 18
Set-TextValue "D18" "0.0000151"
$ws.Range("E18").Value = "  -7.95%  "

# Row 19
Set-TextValue "D19" "393.91"
$ws.Range("E19").Value = "  -8.43%  "

# Row 20
Set-TextValue "D20" "12.28"
$ws.Range("E20").Value = "  -6.01%  "

# Row 21
Set-TextValue "D21" "5.09"
$ws.Range("E21").Value = "  -7.13%  "

# Row 22
Set-TextValue "D22" "7.01"
$ws.Range("E22").Value = "  -4.16%  "

# Row 23
$ws.Range("E23").Value = "  -0.34%  "

# Row 24
Set-TextValue "D24" "5.69"
$ws.Range("E24").Value = "  -0.93%  "

# Row 25
Set-TextValue "D25" "67.19"
$ws.Range("E25").Value = "  -6.02%  "

# Row 26
$ws.Range("E26").Value = "  -5.49%  "

# Row 27
Set-TextValue "D27" "0.475"
$ws.Range("E27").Value = "  -7.73%  "

# Row 28
Set-TextValue "D28" "0.0₃0996"
$ws.Range("E28").Value = "  -13.39%  "

# Row 29 (swapped with what was row 30's content)
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D29" "1.00"
$ws.Range("E29").Value = "  +0.15%  "

# Row 30
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D30" "8.53"
$ws.Range("E30").Value = "  -7.35%  "

# Row 31
$ws.Range("E31").Value = "  -0.11%  "

# Row 32
$ws.Range("E32").Value = "  -8.15%  "

# Row 33
Set-TextValue "D33" "20.75"
$ws.Range("E33").Value = "  -7.35%  "

# Row 34
Set-TextValue "D34" "6.17"
$ws.Range("E34").Value = "  -6.28%  "

# Row 35
Set-TextValue "D35" "4.72"
$ws.Range("E35").Value = "  -8.90%  "

# Row 36
Set-TextValue "D36" "153.51"
$ws.Range("E36").Value = "  -3.39%  "

# Row 37
Set-TextValue "D37" "1.08"
$ws.Range("E37").Value = "  -8.71%  "

# Row 38
Set-TextValue "D38" "1.30"
$ws.Range("E38").Value = "  -9.17%  "

# Row 39
Set-TextValue "D39" "2.698.98"
$ws.Range("E39").Value = "  -6.00%  "

# Row 40
Set-TextValue "D40" "1.63"
$ws.Range("E40").Value = "  -8.74%  "

# Row 41
Set-TextValue "D41" "23.01"
$ws.Range("E41").Value = "  -12.85%  "

# Row 42 (swapped with what was row 43's content)
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D42" "38.08"
$ws.Range("E42").Value = "  -4.40%  "

# Row 43
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D43" "3.97"
$ws.Range("E43").Value = "  -8.45%  "

# Row 44
Set-TextValue "D44" "0.687"
$ws.Range("E44").Value = "  -8.75%  "

# Row 45
Set-TextValue "D45" "0.0597"
$ws.Range("E45").Value = "  -6.23%  "

# Row 46
Set-TextValue "D46" "5.20"
$ws.Range("E46").Value = "  -12.22%  "

# Row 47
Set-TextValue "D47" "0.0252"
$ws.Range("E47").Value = "  -6.79%  "

# Row 48 (swapped with what was row 49's content)
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D48" "1.00"
$ws.Range("E48").Value = "  +0.07%  "

# Row 49
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D49" "20.56"
$ws.Range("E49").Value = "  -10.10%  "

# Row 50
Set-TextValue "D50" "278.05"
$ws.Range("E50").Value = "  -11.12%  "

# Row 51
Set-TextValue "D51" "0.0966"
$ws.Range("E51").Value = "  -5.73%  "
